# Auto-generated Excel COM-interop edit script
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed market data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 158.5
$ws.Range("I6").Value = 113
$ws.Range("K6").Value = 339
$ws.Range("M6").Value = -227
$ws.Range("H11").Value = 60
$ws.Range("I11").Value = 60
$ws.Range("K11").Value = 60
$ws.Range("M11").Value = 80
$ws.Range("H38").Value = 26851.5
$ws.Range("I38").Value = 33468.668
$ws.Range("K38").Value = 100406.004
$ws.Range("M38").Value = -100034.004
$ws.Range("H80").Value = 393.41177
$ws.Range("I80").Value = 431.07693
$ws.Range("K80").Value = 1293.23079
$ws.Range("M80").Value = -295.2307900000001
$ws.Range("H83").Value = 393.41177
$ws.Range("I83").Value = 431.07693
$ws.Range("K83").Value = 3879.69237
$ws.Range("M83").Value = 1112.30763
$ws.Range("H113").Value = 3290.125
$ws.Range("I113").Value = 3180.6
$ws.Range("J113").Value = 3472.6667
$ws.Range("K113").Value = 3180.6
$ws.Range("L113").Value = 3472.6667
$ws.Range("M113").Value = 73.40000000000009
$ws.Range("N113").Value = -9980.6667
$ws.Range("H137").Value = 1194012.6
$ws.Range("I137").Value = 8335665.5
$ws.Range("K137").Value = 25006996.5
$ws.Range("M137").Value = -25004446.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5442.3125
$ws.Range("I32").Value = 5308.96
$ws.Range("J32").Value = 5918.5713
$ws.Range("K32").Value = 5308.96
$ws.Range("L32").Value = 5918.5713
$ws.Range("M32").Value = -5021.96
$ws.Range("N32").Value = -6492.5713
$ws.Range("H61").Value = 3366.6667
$ws.Range("I61").Value = 880
$ws.Range("K61").Value = 880
$ws.Range("M61").Value = -668
$ws.Range("H74").Value = 208227.6
$ws.Range("I74").Value = 618302.5600000001
$ws.Range("K74").Value = 618302.5600000001
$ws.Range("M74").Value = -617428.5600000001
$ws.Range("H77").Value = 208227.6
$ws.Range("I77").Value = 618302.5600000001
$ws.Range("K77").Value = 3091512.8
$ws.Range("M77").Value = -3087144.8
$ws.Range("H122").Value = 1573.125
$ws.Range("I122").Value = 1329.8182
$ws.Range("K122").Value = 3989.4546
$ws.Range("M122").Value = -1539.4546
$ws.Range("H132").Value = 1361.7715
$ws.Range("J132").Value = 4142.857
$ws.Range("L132").Value = 12428.571
$ws.Range("N132").Value = -17488.571
$ws.Range("H136").Value = 3366.6667
$ws.Range("I136").Value = 880
$ws.Range("K136").Value = 2640
$ws.Range("M136").Value = -90

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33339252
$ws.Range("I20").Value = 71437590
$ws.Range("K20").Value = 71437590
$ws.Range("M20").Value = -71437343
$ws.Range("H36").Value = 2562.3333
$ws.Range("I36").Value = 2562.3333
$ws.Range("K36").Value = 2562.3333
$ws.Range("M36").Value = -2028.3333
$ws.Range("H80").Value = 966.6667
$ws.Range("J80").Value = 997
$ws.Range("L80").Value = 997
$ws.Range("N80").Value = -2993
$ws.Range("H83").Value = 966.6667
$ws.Range("J83").Value = 997
$ws.Range("L83").Value = 4985
$ws.Range("N83").Value = -14969
$ws.Range("H86").Value = 2009.1364
$ws.Range("I86").Value = 1677.8572
$ws.Range("K86").Value = 1677.8572
$ws.Range("M86").Value = -554.8571999999999
$ws.Range("H89").Value = 2009.1364
$ws.Range("I89").Value = 1677.8572
$ws.Range("K89").Value = 8389.286
$ws.Range("M89").Value = -2773.286
$ws.Range("H135").Value = 101702.055
$ws.Range("J135").Value = 101702.055
$ws.Range("L135").Value = 101702.055
$ws.Range("N135").Value = -111842.055

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1170
$ws.Range("I19").Value = 715.6667
$ws.Range("J19").Value = 3896
$ws.Range("K19").Value = 715.6667
$ws.Range("L19").Value = 3896
$ws.Range("M19").Value = -545.6667
$ws.Range("N19").Value = -4236
$ws.Range("H24").Value = 1170
$ws.Range("I24").Value = 715.6667
$ws.Range("J24").Value = 3896
$ws.Range("K24").Value = 715.6667
$ws.Range("L24").Value = 3896
$ws.Range("M24").Value = -545.6667
$ws.Range("N24").Value = -4236
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 5555.4546
$ws.Range("I31").Value = 3456
$ws.Range("J31").Value = 8074.8
$ws.Range("K31").Value = 3456
$ws.Range("L31").Value = 8074.8
$ws.Range("M31").Value = -3161
$ws.Range("N31").Value = -8664.799999999999
$ws.Range("H34").Value = 5555.4546
$ws.Range("I34").Value = 3456
$ws.Range("J34").Value = 8074.8
$ws.Range("K34").Value = 3456
$ws.Range("L34").Value = 8074.8
$ws.Range("M34").Value = -3254
$ws.Range("N34").Value = -8478.799999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5139.439
$ws.Range("I131").Value = 13237.833
$ws.Range("J131").Value = 1788.3793
$ws.Range("K131").Value = 39713.499
$ws.Range("L131").Value = 5365.1379
$ws.Range("M131").Value = -34673.499
$ws.Range("N131").Value = -15445.1379

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1348.25
$ws.Range("I102").Value = 1165.8572
$ws.Range("K102").Value = 1165.8572
$ws.Range("M102").Value = 456.1428000000001
$ws.Range("H132").Value = 1873.1852
$ws.Range("I132").Value = 1396.5555
$ws.Range("J132").Value = 2826.4443
$ws.Range("K132").Value = 4189.666499999999
$ws.Range("L132").Value = 8479.332900000001
$ws.Range("M132").Value = -1659.666499999999
$ws.Range("N132").Value = -13539.3329

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7000.6
$ws.Range("I7").Value = 5001
$ws.Range("K7").Value = 5001
$ws.Range("M7").Value = -4889
$ws.Range("H40").Value = 4070.1614
$ws.Range("I40").Value = 4027.9062
$ws.Range("K40").Value = 4027.9062
$ws.Range("M40").Value = -3891.9062
$ws.Range("H68").Value = 3977.8572
$ws.Range("I68").Value = 4411.75
$ws.Range("J68").Value = 3399.3333
$ws.Range("K68").Value = 4411.75
$ws.Range("L68").Value = 3399.3333
$ws.Range("M68").Value = -3662.75
$ws.Range("N68").Value = -4897.3333
$ws.Range("H71").Value = 3977.8572
$ws.Range("I71").Value = 4411.75
$ws.Range("J71").Value = 3399.3333
$ws.Range("K71").Value = 22058.75
$ws.Range("L71").Value = 16996.6665
$ws.Range("M71").Value = -18314.75
$ws.Range("N71").Value = -24484.6665
$ws.Range("H100").Value = 2654.2727
$ws.Range("I100").Value = 1399.8572
$ws.Range("J100").Value = 4849.5
$ws.Range("K100").Value = 1399.8572
$ws.Range("L100").Value = 4849.5
$ws.Range("M100").Value = -858.8571999999999
$ws.Range("N100").Value = -5931.5
$ws.Range("H126").Value = 7000.6
$ws.Range("I126").Value = 5001
$ws.Range("K126").Value = 15003
$ws.Range("M126").Value = -12533
$ws.Range("H132").Value = 5053.778
$ws.Range("I132").Value = 2884
$ws.Range("K132").Value = 8652
$ws.Range("M132").Value = -6122
$ws.Range("H136").Value = 4094.4
$ws.Range("I136").Value = 4126.067
$ws.Range("K136").Value = 12378.201
$ws.Range("M136").Value = -9828.201000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1801.2307
$ws.Range("I126").Value = 1083.2727
$ws.Range("J126").Value = 5750
$ws.Range("K126").Value = 3249.8181
$ws.Range("L126").Value = 17250
$ws.Range("M126").Value = -779.8181
$ws.Range("N126").Value = -22190
$ws.Range("H132").Value = 3410.125
$ws.Range("I132").Value = 4681.75
$ws.Range("J132").Value = 2138.5
$ws.Range("K132").Value = 14045.25
$ws.Range("L132").Value = 6415.5
$ws.Range("M132").Value = -11515.25
$ws.Range("N132").Value = -11475.5
